$d = $word.ActiveDocument

function Get-ParaByText($searchText) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $searchText) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------
# 1. Delete the "_GoBack" bookmark that currently sits inside the
#    Testability paragraph (it will be re-created later near "Rules").
# ---------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------
# 2. Merge the "Tes" / "tability" runs into a single "Testability" run.
#    Insert a throw-away character at the run boundary, then delete a
#    range spanning it so Word coalesces the two identically-formatted
#    runs into one (preserving rFonts/sz/szCs).
# ---------------------------------------------------------------
$pTest = Get-ParaByText("Testability")
$boundary = $pTest.Range.Start + 3   # after "Tes"
$ins = $d.Range($boundary, $boundary)
$ins.InsertAfter("X")
$mergeRng = $d.Range($boundary, $boundary + 1)
$mergeRng.Delete()

# ---------------------------------------------------------------
# 3. Delete the "Hints:" paragraph and the enumeration-hint bullet
#    that follows it.
# ---------------------------------------------------------------
$pHints = Get-ParaByText("Hints:")
$pEnum = $pHints.Next()
$delRng = $d.Range($pHints.Range.Start, $pEnum.Range.End)
$delRng.Delete()

# ---------------------------------------------------------------
# 4. Rule bullet list edits.
# ---------------------------------------------------------------
$pSwitch = Get-ParaByText("Add ability to switch between morning and night and have that be the first, required parameter")
$s = $pSwitch.Range.Start
$e = $pSwitch.Range.End - 1
$r = $d.Range($s, $e)
$r.Text = "Add ability to have different dishes in the morning and at night (See sample input/output below)"

$pFix = Get-ParaByText("Fix existing tests")
$s2 = $pFix.Range.Start
$e2 = $pFix.Range.End - 1
$r2 = $d.Range($s2, $e2)
$fullText = "Make program backwards compatible or fix existing tests"
$r2.Text = $fullText

$part1 = "Make program backwards comp"
$part2 = "atible"
$part3 = " or fix existing tests"

$p2Start = $s2 + $part1.Length
$p2End = $p2Start + $part2.Length
$rSplit2 = $d.Range($p2Start, $p2End)
$rSplit2.Bold = 1
$rSplit2.Bold = 0

$p3Start = $p2End
$p3End = $p3Start + $part3.Length
$rSplit3 = $d.Range($p3Start, $p3End)
$rSplit3.Bold = 1
$rSplit3.Bold = 0

# Delete the now-duplicate "Add ability to have different dishes..." bullet
$pDup = Get-ParaByText("Add ability to have different dishes in the morning and at night (See sample input/output below)")
# There are now two paragraphs with this text (the one we just edited, and
# the original one further down) -- Get-ParaByText returns the first match,
# which is the one we just edited; find the *next* one instead.
$pDupNext = $pDup.Next()
while ($pDupNext -ne $null -and $pDupNext.Range.Text.TrimEnd([char]13) -ne "Add ability to have different dishes in the morning and at night (See sample input/output below)") {
    $pDupNext = $pDupNext.Next()
}
$delRng2 = $d.Range($pDupNext.Range.Start, $pDupNext.Range.End)
$delRng2.Delete()

# ---------------------------------------------------------------
# 5. Re-insert the "_GoBack" bookmark just before the "Rules" run.
# ---------------------------------------------------------------
$pComplete = Get-ParaByText("Complete coverage of the morning dish test cases below")
$pRules = $pComplete.Next()
$bmPos = $pRules.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null

# ---------------------------------------------------------------
# 6. Add a new row ("4 (dessert)" / "cake") to the Morning dishes table.
# ---------------------------------------------------------------
$morningTable = $d.Tables.Item(1)
$newRow = $morningTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "4 (dessert)"
$newRow.Cells.Item(2).Range.Text = "cake"

# ---------------------------------------------------------------
# 7. Move <w:lastRenderedPageBreak/> from the "morning, " run to the
#    "morning, 1, 2, 2" run.
# ---------------------------------------------------------------
Write-Output "done"
